# Apply the "control of nrows" / "ABET" changes to the Service data workbook.
#
# 1. Notes sheet: split the combined "X / Y / Z" option lists that lived in a
#    single B-column cell into separate cells across columns B..E (one value
#    per column), and add a new "Winter" term option.
# 2. Data sheet: add data validation to the new/updated columns -
#    - column B (Type)   -> list of values from Notes!B11:E11
#    - column C (Position)-> list of values from Notes!B12:C12
#    - column D (Term)    -> list of values from Notes!B13:E13
#    - column E (Calendar Year) -> whole number between 1900 and 2100
#    - column F (Hours/Semester) -> decimal between 0 and 100000

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Notes sheet - split the combined option strings into individual cells
# ---------------------------------------------------------------------
$notes = $wb.Worksheets.Item("Notes")

# Row 11 - Type: "Professional, University, Department, Community"
$notes.Range("B11").Value2 = "Professional"
$notes.Range("C11").Value2 = "University"
$notes.Range("D11").Value2 = "Department"
$notes.Range("E11").Value2 = "Community"

# Row 12 - Position: "Chair / Member"
$notes.Range("B12").Value2 = "Chair"
$notes.Range("C12").Value2 = "Member"

# Row 13 - Term: "Fall / Spring / Summer" (+ new "Winter" option)
$notes.Range("B13").Value2 = "Fall"
$notes.Range("C13").Value2 = "Spring"
$notes.Range("D13").Value2 = "Summer"
$notes.Range("E13").Value2 = "Winter"

$notes.Range("B13").Select() | Out-Null

# ---------------------------------------------------------------------
# Data sheet - wire up data validation against the Notes option lists
# and add numeric range validation for Calendar Year / Hours per Semester
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

$data.Range("B1:B1048576").Validation.Add(3, 1, 1, "=Notes!`$B`$11:`$E`$11")
$data.Range("C1:C1048576").Validation.Add(3, 1, 1, "=Notes!`$B`$12:`$C`$12")
$data.Range("D1:D1048576").Validation.Add(3, 1, 1, "=Notes!`$B`$13:`$E`$13")
$data.Range("E1:E1048576").Validation.Add(1, 1, 1, 1900, 2100)
$data.Range("F1:F1048576").Validation.Add(2, 1, 1, 0, 100000)

$data.Range("E6").Select() | Out-Null
